$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.99999999361158021
$ws.Range("A2").Value = 0.99797072345550553
$ws.Range("A3").Value = 0.99330776872850968
$ws.Range("A4").Value = 0.99659562760113241
$ws.Range("A5").Value = 0.98674669336091059
$ws.Range("A6").Value = 0.96374123627344188
$ws.Range("A7").Value = 0.96508416963765664
$ws.Range("A8").Value = 0.96561983692837949
$ws.Range("A9").Value = 0.96241650340408358
$ws.Range("A10").Value = 0.96014245088277539
$ws.Range("A11").Value = 0.95986021509951258
$ws.Range("A12").Value = 0.95957931299260968
$ws.Range("A13").Value = 0.96066160114794941
$ws.Range("A14").Value = 0.96215231033011883
$ws.Range("A15").Value = 0.96475329535882726
$ws.Range("A16").Value = 0.96834341745612973
$ws.Range("A17").Value = 0.96463593087383337
$ws.Range("A18").Value = 0.96352704148823265
$ws.Range("A19").Value = 0.99691751902046999
$ws.Range("A20").Value = 0.98980055316335513
$ws.Range("A21").Value = 0.98840207433825755
$ws.Range("A22").Value = 0.9871375751627256
$ws.Range("A23").Value = 0.98491469600069448
$ws.Range("A24").Value = 0.96883989120913805
$ws.Range("A25").Value = 0.95939419948421945
$ws.Range("A26").Value = 0.96008812432044666
$ws.Range("A27").Value = 0.9552445854616094
$ws.Range("A28").Value = 0.93378164092238047
$ws.Range("A29").Value = 0.91851391824175976
$ws.Range("A30").Value = 0.91194427691334901
$ws.Range("A31").Value = 0.90429107819952403
$ws.Range("A32").Value = 0.90261178773893525
$ws.Range("A33").Value = 0.90209178626789199
